$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: 展览 (Exhibition) - "want to go" counter (column F) bumped on
# several rows. No structural change.
# ---------------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item(1)
$wsExpo.Cells.Item(5,6).Value  = 8699
$wsExpo.Cells.Item(7,6).Value  = 11040
$wsExpo.Cells.Item(9,6).Value  = 16
$wsExpo.Cells.Item(22,6).Value = 1870
$wsExpo.Cells.Item(23,6).Value = 698
$wsExpo.Cells.Item(30,6).Value = 1262
$wsExpo.Cells.Item(32,6).Value = 8
$wsExpo.Cells.Item(36,6).Value = 456
$wsExpo.Cells.Item(38,6).Value = 297
$wsExpo.Cells.Item(42,6).Value = 371
$wsExpo.Cells.Item(43,6).Value = 111

# ---------------------------------------------------------------------------
# Sheet 2: 演出 (Performance) - a new concert ("黄霄雲" 2024 tour) was added
# to the source feed; because the list is sorted by start date and this show
# starts 2024-11-01 it lands right before the existing 2024-11-03 entry, i.e.
# as the new row 19. Every row from the old row 19 onward shifts down by one.
# Column F14 also got its counter bumped.
# ---------------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item(2)

$wsShow.Cells.Item(14,6).Value = 26

# Insert a blank row at position 19, pushing old rows 19-23 down to 20-24.
$wsShow.Rows.Item(19).Insert()

# Copy column-A's number formatting/border/bold style from the row below
# (which holds the old row-19 formatting) onto the freshly inserted row so
# the new id cell matches the rest of the table.
$wsShow.Cells.Item(20,1).Copy()
$wsShow.Cells.Item(19,1).PasteSpecial(-4122)

# Fill in the new row 19.
$wsShow.Cells.Item(19,1).Value = 18
$wsShow.Cells.Item(19,2).NumberFormat = "@"
$wsShow.Cells.Item(19,2).Value = "2024-11-01"
$wsShow.Cells.Item(19,2).Style = "Normal"
$wsShow.Cells.Item(19,3).Value = "北京·黄霄雲「宇宙无敌号」2024 演唱会"
$wsShow.Cells.Item(19,4).Value = "复兴路69号(原五棵松体育馆) 华熙LIVE凯迪拉克中心"
$wsShow.Cells.Item(19,5).Value = "2024.11.01 19:30-11.02 22:00"
$wsShow.Cells.Item(19,6).Value = 45
$wsShow.Cells.Item(19,7).Value = 388
$wsShow.Cells.Item(19,8).Value = "https://show.bilibili.com/platform/detail.html?id=92913"
$wsShow.Cells.Item(19,9).Value = "//i0.hdslb.com/bfs/openplatform/202409/mnAmm0Wu1727257608727.jpeg"

# The plain row-insert keeps each shifted row's own column-A id value as-is
# (it does not renumber). The published sheet instead keeps the sequential
# "row-1" id numbering intact, so re-stamp column A for rows 19-24.
$wsShow.Cells.Item(19,1).Value = 18
$wsShow.Cells.Item(20,1).Value = 19
$wsShow.Cells.Item(21,1).Value = 20
$wsShow.Cells.Item(22,1).Value = 21
$wsShow.Cells.Item(23,1).Value = 22
$wsShow.Cells.Item(24,1).Value = 23

# ---------------------------------------------------------------------------
# Sheet 3: 本地生活 (Local life) - counter bump.
# ---------------------------------------------------------------------------
$wsLocal = $wb.Worksheets.Item(3)
$wsLocal.Cells.Item(3,6).Value = 2833

# ---------------------------------------------------------------------------
# Sheet 4: 全部类型 (All types, merged view) - counter bumps mirroring the
# same underlying events updated on sheet 1 / sheet 2 above.
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Cells.Item(8,6).Value  = 8699
$wsAll.Cells.Item(10,6).Value = 11040
$wsAll.Cells.Item(19,6).Value = 1870
$wsAll.Cells.Item(20,6).Value = 698
$wsAll.Cells.Item(29,6).Value = 1262
$wsAll.Cells.Item(31,6).Value = 8
$wsAll.Cells.Item(33,6).Value = 26
$wsAll.Cells.Item(38,6).Value = 456
$wsAll.Cells.Item(42,6).Value = 371
$wsAll.Cells.Item(43,6).Value = 111
